$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while keeping the cell a plain TEXT cell (no residual
# style/number-format change) — mirrors the source data's inline-string cells.
# Numeric-looking strings (e.g. "18.17", "0.787") would otherwise be silently
# coerced to actual numbers by the normal .Value setter, which both changes
# the cell type and can introduce floating point noise. Forcing the cell to
# the "Text" number format before the write prevents that coercion, and
# resetting the style back to "Normal" afterwards removes the temporary
# formatting change so the cell ends up identical in style to before.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "26.260.29"
Set-TextValue "E2" "  +2.79%  "
Set-TextValue "D3" "1.609.66"
Set-TextValue "E3" "  +1.32%  "
Set-TextValue "E4" "  -0.58%  "
Set-TextValue "D5" "212.91"
Set-TextValue "E5" "  +2.08%  "
Set-TextValue "E6" "  -0.66%  "
Set-TextValue "E7" "  +0.78%  "
Set-TextValue "E8" "  +1.66%  "
Set-TextValue "D9" "0.0619"
Set-TextValue "E9" "  +1.63%  "
Set-TextValue "D10" "18.17"
Set-TextValue "E10" "  +2.53%  "
Set-TextValue "D11" "0.0815"
Set-TextValue "E11" "  +4.17%  "
Set-TextValue "D12" "1.833.76"
Set-TextValue "E12" "  +1.25%  "
Set-TextValue "D13" "1.608.81"
Set-TextValue "E13" "  +1.29%  "
Set-TextValue "D14" "4.02"
Set-TextValue "E14" "  -0.46%  "
Set-TextValue "D15" "0.510"
Set-TextValue "E15" "  +0.61%  "
Set-TextValue "D16" "26.224.49"
Set-TextValue "E16" "  +2.49%  "
Set-TextValue "D17" "60.72"
Set-TextValue "E17" "  +0.68%  "
Set-TextValue "D18" "0.0₃0729"
Set-TextValue "E18" "  +2.63%  "
Set-TextValue "E19" "  -0.45%  "
Set-TextValue "D20" "199.25"
Set-TextValue "E20" "  +6.19%  "
Set-TextValue "E21" "  +2.09%  "
Set-TextValue "D22" "9.41"
Set-TextValue "E22" "  +1.00%  "
Set-TextValue "D23" "6.03"
Set-TextValue "E23" "  +1.76%  "
Set-TextValue "D24" "0.131"
Set-TextValue "E24" "  +2.10%  "
Set-TextValue "D25" "142.64"
Set-TextValue "E25" "  +1.50%  "
Set-TextValue "D26" "1.76"
Set-TextValue "E26" "  +4.02%  "
Set-TextValue "E27" "  -0.52%  "
Set-TextValue "D28" "15.20"
Set-TextValue "E28" "  +1.88%  "
Set-TextValue "E29" "  +0.05%  "
Set-TextValue "E30" "  -0.71%  "
Set-TextValue "E31" "  +1.46%  "
Set-TextValue "D32" "3.14"
Set-TextValue "E32" "  +2.54%  "
Set-TextValue "E33" "  +0.57%  "
Set-TextValue "E34" "  +1.84%  "
Set-TextValue "E35" "  -1.54%  "
Set-TextValue "D36" "1.108.07"
Set-TextValue "E36" "  +1.85%  "
Set-TextValue "E37" "  -0.37%  "
Set-TextValue "E38" "  +1.11%  "
Set-TextValue "E39" "  -0.59%  "
Set-TextValue "B40" "ARBITRUM"
Set-TextValue "C40" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D40" "0.787"
Set-TextValue "E40" "  +1.21%  "
Set-TextValue "B41" "ImmutableX"
Set-TextValue "C41" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D41" "0.501"
Set-TextValue "E41" "  +1.46%  "
Set-TextValue "D42" "0.783"
Set-TextValue "E42" "  +6.74%  "
Set-TextValue "D43" "1.744.59"
Set-TextValue "E43" "  +1.18%  "
Set-TextValue "E44" "  +0.99%  "
Set-TextValue "D45" "92.59"
Set-TextValue "E45" "  -2.72%  "
Set-TextValue "E46" "  +1.35%  "
Set-TextValue "E47" "  +9.18%  "
Set-TextValue "E48" "  +1.68%  "
Set-TextValue "E49" "  +0.16%  "
Set-TextValue "E50" "  +0.33%  "
Set-TextValue "E51" "  -0.33%  "
